# Auto commit edit: refresh the raw metric totals on the "Metrics" sheet.
# Everything on the "today" sheet (B11:B22, E11:E22, F11:F22) is formula
# driven off these cells (directly or via Metrics!B2:B13), and A1's
# "=TODAY()-1" is volatile, so they recompute on their own once the
# workbook recalculates after this script runs - no need to touch them
# by hand.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# New raw values for Metrics!B2:B13 (same order as the sheet).
$wsMetrics.Range("B2").Value  = 152769.86000000002
$wsMetrics.Range("B3").Value  = 115948.85
$wsMetrics.Range("B4").Value  = 37448.200000000004
$wsMetrics.Range("B5").Value  = 6171
$wsMetrics.Range("B6").Value  = 5788640.5899999999
$wsMetrics.Range("B7").Value  = 4886666.4800000004
$wsMetrics.Range("B8").Value  = 1701540.02
$wsMetrics.Range("B9").Value  = 226448
$wsMetrics.Range("B10").Value = 34254021.579999998
$wsMetrics.Range("B11").Value = 32161941.640000004
$wsMetrics.Range("B12").Value = 11983262.059999999
$wsMetrics.Range("B13").Value = 1324078

# Restore the cursor/selection on Metrics first ...
$wsMetrics.Range("D15").Select() | Out-Null

# ... then select on "today" last, so it stays the active sheet/tab
# (matches the workbook's original activeTab / tabSelected="1").
$wsToday.Range("E8").Select() | Out-Null
